$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 1998
$ws.Range("E2").Value = -85
$ws.Range("F2").Value = -85
$ws.Range("G2").Value = -132
$ws.Range("H2").Value = -124
$ws.Range("I2").Value = -124
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 1858
$ws.Range("L2").Value = 993
$ws.Range("M2").Value = 865
$ws.Range("N2").Value = 865
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 407
$ws.Range("Q2").Value = -138
$ws.Range("R2").Value = -16
$ws.Range("S2").Value = 27
$ws.Range("T2").Value = 19
$ws.Range("U2").Value = -158
$ws.Range("V2").Value = 631
$ws.Range("W2").Value = -4.25
$ws.Range("X2").Value = -6.2
$ws.Range("Y2").Value = -13.33
$ws.Range("Z2").Value = -6.16
$ws.Range("AA2").Value = 114.77
$ws.Range("AB2").Value = 112.4
$ws.Range("AC2").Value = -137
$ws.Range("AD2").Value = -7.14
$ws.Range("AE2").Value = 955
$ws.Range("AF2").Value = 1.02
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 90569004

# --- Row 3 ---
$ws.Range("D3").Value = 1713
$ws.Range("E3").Value = 32
$ws.Range("F3").Value = 32
$ws.Range("G3").Value = 24
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 23
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 1775
$ws.Range("L3").Value = 891
$ws.Range("M3").Value = 884
$ws.Range("N3").Value = 884
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 407
$ws.Range("Q3").Value = 43
$ws.Range("R3").Value = 12
$ws.Range("S3").Value = -11
$ws.Range("T3").Value = 18
$ws.Range("U3").Value = 25
$ws.Range("V3").Value = 557
$ws.Range("W3").Value = 1.9
$ws.Range("X3").Value = 1.35
$ws.Range("Y3").Value = 2.65
$ws.Range("Z3").Value = 1.28
$ws.Range("AA3").Value = 100.73
$ws.Range("AB3").Value = 117.07
$ws.Range("AC3").Value = 26
$ws.Range("AD3").Value = 48.06
$ws.Range("AE3").Value = 976
$ws.Range("AF3").Value = 1.26
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 90569004

# --- Row 4 ---
$ws.Range("D4").Value = 1293
$ws.Range("E4").Value = -120
$ws.Range("F4").Value = -120
$ws.Range("G4").Value = -146
$ws.Range("H4").Value = -147
$ws.Range("I4").Value = -147
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 1598
$ws.Range("L4").Value = 792
$ws.Range("M4").Value = 805
$ws.Range("N4").Value = 805
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 428
$ws.Range("Q4").Value = 68
$ws.Range("R4").Value = -9
$ws.Range("S4").Value = -84
$ws.Range("T4").Value = 18
$ws.Range("U4").Value = 50
$ws.Range("V4").Value = 437
$ws.Range("W4").Value = -9.26
$ws.Range("X4").Value = -11.33
$ws.Range("Y4").Value = -17.34
$ws.Range("Z4").Value = -8.69
$ws.Range("AA4").Value = 98.39
$ws.Range("AB4").Value = 88.18000000000001
$ws.Range("AC4").Value = -159
$ws.Range("AD4").Value = -7.24
$ws.Range("AE4").Value = 846
$ws.Range("AF4").Value = 1.36
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 95163761

# --- Row 5 ---
$ws.Range("D5").Value = 1471
$ws.Range("E5").Value = -65
$ws.Range("F5").Value = -65
$ws.Range("G5").Value = -308
$ws.Range("H5").Value = -306
$ws.Range("I5").Value = -306
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 1504
$ws.Range("L5").Value = 974
$ws.Range("M5").Value = 530
$ws.Range("N5").Value = 530
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 436
$ws.Range("Q5").Value = -46
$ws.Range("R5").Value = 9
$ws.Range("S5").Value = 59
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = -50
$ws.Range("V5").Value = 482
$ws.Range("W5").Value = -4.41
$ws.Range("X5").Value = -20.8
$ws.Range("Y5").Value = -45.85
$ws.Range("Z5").Value = -19.74
$ws.Range("AA5").Value = 183.91
$ws.Range("AB5").Value = 21.43
$ws.Range("AC5").Value = -318
$ws.Range("AD5").Value = -2.36
$ws.Range("AE5").Value = 546
$ws.Range("AF5").Value = 1.38
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 97001667

# --- Row 6 ---
$ws.Range("D6").Value = 1620
$ws.Range("E6").Value = -40
$ws.Range("F6").Value = -40
$ws.Range("G6").Value = -66
$ws.Range("H6").Value = -66
$ws.Range("I6").Value = -66
$ws.Range("K6").Value = 1636
$ws.Range("L6").Value = 846
$ws.Range("M6").Value = 791
$ws.Range("N6").Value = 791
$ws.Range("P6").Value = 595
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = -115
$ws.Range("S6").Value = 233
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 371
$ws.Range("W6").Value = -2.46
$ws.Range("X6").Value = -4.06
$ws.Range("Y6").Value = -9.960000000000001
$ws.Range("Z6").Value = -4.19
$ws.Range("AA6").Value = 106.97
$ws.Range("AB6").Value = 33
$ws.Range("AC6").Value = -66
$ws.Range("AD6").Value = -22.97
$ws.Range("AE6").Value = 665
$ws.Range("AF6").Value = 2.29
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 118885290

# --- Rows 7-9: clear all data columns (D..AI), keep A/B/C and AJ as-is (AJ also cleared since it had no value before) ---
$ws.Range("D7:AI9").ClearContents()
